# Fix the corrupted label in Tabelle2!A2.
# It currently holds the garbled shared string "Liq2342332423uide Mittel"
# and should instead read "Liquide Mittel" (matching Tabelle1!A2 and the
# other clean "Liquide Mittel" entry already present in the shared string
# table). Writing the clean text here lets Excel reuse the existing shared
# string and drop the now-unreferenced garbled one on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")
$ws.Range("A2").Value = "Liquide Mittel"
